$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting B:F (old A:E) to the right
$ws.Columns.Item(1).Insert()

# Header for new ID column, matching the style of the other header cells
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)  # xlPasteFormats

# Row labels for the ID column
$labels = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
